$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actividades")
$ws.Range("D16").Value = "test"
